# Regenerate save_data: update column G ("K") values (was "Strike#"-derived), rows 2-25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 3
    4  = 2
    5  = 3
    6  = 4
    7  = 2
    8  = 1
    9  = 5
    10 = 6
    11 = 1
    12 = 3
    13 = 2
    14 = 5
    15 = 3
    16 = 5
    17 = 2
    18 = 0
    19 = 6
    20 = 2
    21 = 3
    22 = 3
    23 = 3
    24 = 4
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
